$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------
# Overview sheet: per-language status columns (E = zh-cn, F = de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-language detail sheets: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the status columns -----------------------------------------------
# Target stored width (OOXML <col width="...">) is 13.4101845877511, down from
# 17.2159881591797. ColumnWidth is quantized on export to 1/6-character steps
# (stored width = round(ColumnWidth*6)/6 + 5/6), so feed the pre-offset value
# whose rounded result lands closest to the target width.
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # F (de-de status)

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # C (Status)

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # C (Status)
